# Update Name of Algo
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E6").Value = 16.316
$ws.Range("D7").Value = -7.198
$ws.Range("B9").Value = 5.859
$ws.Range("D12").Value = -7.269
$ws.Range("D14").Value = -7.505
$ws.Range("E15").Value = 16.158
$ws.Range("B18").Value = 5.166
$ws.Range("B20").Value = 6.978
$ws.Range("D26").Value = -7.633
$ws.Range("B27").Value = 6.2
$ws.Range("D27").Value = -8.149000000000001
$ws.Range("D29").Value = -7.348999999999999
$ws.Range("E33").Value = 17.079
$ws.Range("B35").Value = 8.170999999999999
$ws.Range("E35").Value = 16.461
$ws.Range("D37").Value = -7.822
$ws.Range("D38").Value = -7.228
$ws.Range("E38").Value = 16.739
$ws.Range("E43").Value = 17.127
$ws.Range("E44").Value = 17.04199999999999
$ws.Range("E47").Value = 16.459
$ws.Range("D51").Value = -8.401
$ws.Range("E51").Value = 16.634
$ws.Range("D52").Value = -7.831
$ws.Range("D55").Value = -8.129000000000001
$ws.Range("E57").Value = 16.471
$ws.Range("E63").Value = 17.601
$ws.Range("B69").Value = 5.992
$ws.Range("D69").Value = -7.113000000000001
$ws.Range("D70").Value = -7.176
$ws.Range("E70").Value = 17.792
$ws.Range("B76").Value = 6.308
$ws.Range("B78").Value = 8.550999999999998
$ws.Range("D81").Value = -7.505000000000001
$ws.Range("B82").Value = 5.366000000000001
$ws.Range("B83").Value = 5.129
$ws.Range("D83").Value = -8.516999999999999
$ws.Range("E88").Value = 16.385
$ws.Range("B93").Value = 6.209
$ws.Range("E99").Value = 16.768
$ws.Range("D102").Value = -7.865
